$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove existing hyperlink objects so we can recreate them pointing at
#     their (post-insert) shifted rows, in the same order as before so the
#     relationship ids (rId1..rId9) line up the same as before.
$ws.Hyperlinks.Delete()

# --- Insert a new row for "Appearance with e-filing" right after the
#     existing "Appearance" row (row 2), pushing everything else down one.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value2 = "Appearance with e-filing"
$ws.Range("B3").Value2 = "https://www.illinoislegalaid.org/legal-information/appearance-e-filing"
# Give B3 the same visual "Hyperlink" look as the other url cells, even
# though (per the source data) it is not wired up as a clickable hyperlink.
$ws.Range("B3").Style = "Hyperlink"

# --- Recreate the hyperlinks on the (now shifted down by one row) cells.
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.illinoislegalaid.org/legal-information/appearance") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.illinoislegalaid.org/legal-information/fee-waiver") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://www.illinoislegalaid.org/legal-information/collection-proof-debtor-letter") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.illinoislegalaid.org/legal-information/request-collection-agency-stop-contacting") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.illinoislegalaid.org/legal-information/end-illegal-lockout-demand") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter") | Out-Null

# Re-applying the "Hyperlink" cell style ensures these cells keep using the
# same style index as before instead of a duplicate one created by Add().
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"
$ws.Range("B10").Style = "Hyperlink"
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("B12").Style = "Hyperlink"

# --- Refresh the remembered sort range/state to cover the new row (the
#     stored sort range has always lagged one row behind the full data
#     range - e.g. before the edit it was A2:B10 while data ran to B11).
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A11"))
$sort.SetRange($ws.Range("A2:B11"))
$sort.Header = 0
$sort.Apply()

# --- Match the saved selection in the edited workbook.
$ws.Range("C3").Select()
